# Learning Data Entry ( 22-7-2025 )
# Adds MAX / MIN / AVG / COUNT helper columns (P:S) to the "Practice"
# worksheet for rows 7-16, updates the "Track" log on the "Track" sheet
# with a new entry describing the work done, and makes "Practice" the
# active sheet/tab (mirroring the saved UI state captured in the diff).

$wb = $excel.ActiveWorkbook
$practice = $wb.Worksheets.Item("Practice")
$track = $wb.Worksheets.Item("Track")

# ---------------------------------------------------------------------
# 1. Practice sheet: header row (row 6) gains MAX/MIN/AVG/COUNT labels
#    in P6:S6, using the same style as the existing SUB1/SUB2/SUB3/TOTAL
#    headers (L6:O6).
# ---------------------------------------------------------------------
$practice.Range("L6").Copy() | Out-Null
$practice.Range("P6:S6").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$practice.Range("P6").Value = "MAX"
$practice.Range("Q6").Value = "MIN"
$practice.Range("R6").Value = "AVG"
$practice.Range("S6").Value = "COUNT"

# ---------------------------------------------------------------------
# 2. Rows 7-16: add the MAX / MIN / AVERAGE / COUNT formulas over the
#    SUB1:SUB3 (L:N) range for that row.
# ---------------------------------------------------------------------
for ($r = 7; $r -le 16; $r++) {
    $practice.Range("P${r}").Formula = "=MAX(L${r}:N${r})"
    $practice.Range("Q${r}").Formula = "=MIN(L${r}:N${r})"
    $practice.Range("R${r}").Formula = "=AVERAGE(L${r}:N${r})"
    $practice.Range("S${r}").Formula = "=COUNT(L${r}:N${r})"
}

# ---------------------------------------------------------------------
# 3. Track sheet: log the day's lesson (row 6, columns A/B/D) - extend
#    the description text in D6 to mention the new topics covered.
# ---------------------------------------------------------------------
$track.Range("D6").Value = "SUM Formula , adding marks ,paste values, fix refrence , relative refrence,MAX,MIN,AVG,COUNT"

# ---------------------------------------------------------------------
# 4. View state: make "Practice" the active/selected sheet (instead of
#    "Track"), matching the saved workbook view captured in the diff.
# ---------------------------------------------------------------------
$practice.Activate()
$practice.Application.ActiveWindow.ScrollRow = 3
$practice.Application.ActiveWindow.ScrollColumn = 3
$practice.Range("M21").Select() | Out-Null
